$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Septiembre de 2020 a las 20:13"

# --- Update country stat rows with refreshed data ---
# Each entry: row number -> @(Country, TotalCases, NewCases, ActiveCases, Recovered, Critical, DeathsToday, Deaths)
# A few countries (Marruecos/Panama, Etiopia/Portugal, Islas Malvinas/Montserrat)
# swapped rank/rows as part of the refresh.
$rows = @{
    4   = @("Estados Unidos", 7159042, 19489, 4413984, 2538008, 0, 457, 207050)
    5   = @("India",          5793482, 63298, 4732300,  969366, 0, 643,  91816)
    11  = @("España",          704209, 10653,       0,       0, 0,  84,  31118)
    35  = @("Marruecos",       110099,  2356,   90186,   17957, 0,  38,   1956)
    36  = @("Panama",          107990,     0,   84437,   21262, 0,   0,   2291)
    51  = @("Etiopia",          71687,   604,   29461,   41078, 0,   7,   1148)
    52  = @("Portugal",         71156,   691,   46676,   22549, 0,   3,   1931)
    73  = @("Irlanda",          33994,   318,   23364,    8833, 0,   3,   1797)
    112 = @("Mozambique",        7399,   137,    4558,    2790, 0,   2,     51)
    115 = @("Malaui",            5747,     1,    4163,    1405, 0,   0,    179)
    215 = @("Islas Malvinas",      13,     0,      13,       0, 0,   0,      0)
    216 = @("Montserrat",          13,     0,      12,       0, 0,   0,      1)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
